# Replace the placeholder "Team X" entries in the "Team Members" column of
# the Sprint / Product Backlog table with the actual member names, in
# document (row) order.
#
# Table 2 on the ActiveDocument is the sprint backlog table:
#   Sprint | Functional Requirement (Epic) | User Story Number |
#   User Story / Task | Story Points | Priority | Team Members
#
# Rows 2-14 (row 1 is the header) hold the USN-1 .. USN-13 user stories, and
# column 7 is "Team Members" which currently contains "Team A/B/C/D/E"
# placeholders that must become the real assignee names.

$d = $word.ActiveDocument

$table = $d.Tables(2)

$newNames = @(
    "Nihalika",  # row 2  - USN-1  (was Team A)
    "Rishav",    # row 3  - USN-2  (was Team A)
    "Asmita",    # row 4  - USN-3  (was Team B)
    "Vedika",    # row 5  - USN-4  (was Team A)
    "Nihalika",  # row 6  - USN-5  (was Team A)
    "Asmita",    # row 7  - USN-6  (was Team C)
    "Rishav",    # row 8  - USN-7  (was Team C)
    "Vedika",    # row 9  - USN-8  (was Team C)
    "Asmita",    # row 10 - USN-9  (was Team C)
    "Vedika",    # row 11 - USN-10 (was Team D)
    "Nihalika",  # row 12 - USN-11 (was Team D)
    "Rishav",    # row 13 - USN-12 (was Team E)
    "Nihalika"   # row 14 - USN-13 (was Team E)
)

for ($i = 0; $i -lt $newNames.Length; $i++) {
    $rowIndex = $i + 2
    $cell = $table.Cell($rowIndex, 7)
    $cellRange = $cell.Range
    $cellRange.MoveEnd(1, -1) | Out-Null
    $cellRange.Text = $newNames[$i]
}
